# Update cryptocurrency price/volume data (scraper refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.030.61'
$ws.Range("E2").Value = '  -0.91%  '

$ws.Range("D3").Value = '1.906.91'
$ws.Range("E3").Value = '  -0.79%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.31%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7594'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.74'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.32%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3085'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.80%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.54'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.60%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06905'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.11%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08018'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.29%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7563'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.31%  '

$ws.Range("D13").Value = '1.903.75'
$ws.Range("E13").Value = '  -1.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.264'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.80%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.88'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.204'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.65%  '

$ws.Range("D17").Value = '30.034.70'
$ws.Range("E17").Value = '  -0.95%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.06'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.47%  '

$ws.Range("E19").Value = '  -1.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.62'
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9998'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.30%  '

$ws.Range("D22").Value = '2.153.47'
$ws.Range("E22").Value = '  -0.68%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.050'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.74%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.325'
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.72%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.84%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1303'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.080'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.351'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.90%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.525'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.58%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.315'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.99%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.052'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.72%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05428'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.42%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.289'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7390'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.66%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.717'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01948'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.13%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.765'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.88%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.256'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.67%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4460'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.94'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.13%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.951'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9999'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.24%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8318'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.85%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.700'
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.68'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.07%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.864'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.40%  '

$ws.Range("D49").Value = '2.056.83'
$ws.Range("E49").Value = '  -0.97%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.57'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.85%  '

$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '923.50'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.95%  '
